# Apply the change described by the diff:
# Insert a new data row at row 388 in Sheet1 (pushing existing rows 388-486 down
# to 389-487), and populate the newly inserted row 388 with new values while all
# the other columns (A,B,C,E,F,G,H,I,R) keep the same values as neighboring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 388; this shifts rows 388:486 down to 389:487,
# carrying their full contents/formatting with them.
$ws.Rows("388:388").Insert()

# Populate the newly inserted row 388 with the new record's values.
$ws.Cells.Item(388, 1).Value = 5
$ws.Cells.Item(388, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(388, 3).Value = 'Maule'
$ws.Cells.Item(388, 4).Value = 44943
$ws.Cells.Item(388, 5).Value = 7
$ws.Cells.Item(388, 6).Value = 100112032
$ws.Cells.Item(388, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(388, 8).Value = 'Sin especificar'
$ws.Cells.Item(388, 9).Value = 'Primera'
$ws.Cells.Item(388, 10).Value = 300
$ws.Cells.Item(388, 11).Value = 6000
$ws.Cells.Item(388, 12).Value = 6000
$ws.Cells.Item(388, 13).Value = 6000
$ws.Cells.Item(388, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(388, 15).Value = 'Región del Maule'
$ws.Cells.Item(388, 16).Value = 120
$ws.Cells.Item(388, 17).Value = 50
$ws.Cells.Item(388, 18).Value = 'Hortaliza'

# Make sure the D column of the new row keeps the same date-like number
# format ("style 2") as the rest of the D column.
$ws.Range("D388").NumberFormat = $ws.Range("D389").NumberFormat
